# The commit removes the "TC004 / Pharmacy\Reports\Sales\TC01BillWiseSalesReport.py.py"
# test-case row from the "test" worksheet. That row was row 4 in the sheet
# (TestCaseName/Status/Plan/Run No./TC No. header is row 1, data starts at row 2).
# Deleting the entire row shifts all subsequent rows up by one, which matches
# the diff (row 11 disappears, dimension shrinks from A1:I11 to A1:I10, the
# "I6" note cell becomes "I5", and the shared-strings table loses the two
# strings that are no longer referenced anywhere: "TC004" and
# "Pharmacy\Reports\Sales\TC01BillWiseSalesReport.py.py").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the BillWiseSalesReport (TC004) row entirely; everything below
# shifts up automatically.
$ws.Rows.Item(4).Delete()

# The saved selection in the workbook moved to G8 after the edit.
$ws.Range("G8").Select()
